$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching the style of the
# existing header cells (bold, centered, thin border) by copying the
# format from H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I0 and IF values for each data row (2-27).
$data = @(
  @(7,7),
  @(7,8),
  @(6,6),
  @(6,6),
  @(1,1),
  @(7,7),
  @(7,7),
  @(9,9),
  @(6,6),
  @(4,5),
  @(1,2),
  @(9,9),
  @(1,2),
  @(6,7),
  @(6,7),
  @(4,4),
  @(8,8),
  @(1,1),
  @(11,11),
  @(8,8),
  @(5,5),
  @(7,7),
  @(4,4),
  @(8,8),
  @(8,8),
  @(7,7)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
